$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Address each data row by its (unique) table position rather than by searching
# for its text, since several of the original cell contents ("frfrfr" /
# "frfrfrfr", "43434" / "43434 wewe", ...) are substrings of one another and a
# global Find/Replace would clobber neighbouring cells.

$t.Cell(1, 1).Range.Text  = "wq qw wq qw"
$t.Cell(2, 1).Range.Text  = "wq 23/32q"
$t.Cell(3, 1).Range.Text  = "1111111111 eqe"
$t.Cell(4, 1).Range.Text  = "wdw"
$t.Cell(5, 1).Range.Text  = "dwd@dede"
$t.Cell(6, 1).Range.Text  = "+380984343994"
$t.Cell(7, 1).Range.Text  = "dedeed"
$t.Cell(8, 1).Range.Text  = "Amount USD: -610080234"
$t.Cell(9, 1).Range.Text  = "___________________16-4-2020"
$t.Cell(11, 1).Range.Text = "wq qw wq qw"
